# Generate Report for Handoff
# Updates the "b.md" row across the Overview / zh-cn / de-de sheets to
# reflect that the handoff package is now ready, with fresh handoff
# filenames/timestamps and an error detail message for the stale handback.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-08-29 12:40:41"

# ---- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "False"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$ws.Range("H3").Value = "2016-08-29 12:40:36"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30c3429b036f9cabb943d14565e1be0f1b287b5c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19ec0415e550fcd8d22fbd452f32037ad3970cc4/e2e/b.md."
$ws.Columns.Item(16).ColumnWidth = 40

# ---- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("F3").Value = "False"
$ws.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$ws.Range("H3").Value = "2016-08-29 12:40:41"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30c3429b036f9cabb943d14565e1be0f1b287b5c/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19ec0415e550fcd8d22fbd452f32037ad3970cc4/e2e/b.md."
$ws.Columns.Item(16).ColumnWidth = 40
